$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" "42.414.13"
$ws.Range("E2").Value = "  +6.37%  "

# Row 3
Set-TextValue $ws "D3" "2.246.25"
$ws.Range("E3").Value = "  +2.13%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
Set-TextValue $ws "D5" "232.67"
$ws.Range("E5").Value = "  +2.19%  "

# Row 6
Set-TextValue $ws "D6" "0.631"
$ws.Range("E6").Value = "  +0.37%  "

# Row 7
Set-TextValue $ws "D7" "62.21"
$ws.Range("E7").Value = "  -2.15%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
Set-TextValue $ws "D9" "0.406"
$ws.Range("E9").Value = "  +3.15%  "

# Row 10
Set-TextValue $ws "D10" "59.39"
$ws.Range("E10").Value = "  +1.42%  "

# Row 11
Set-TextValue $ws "D11" "0.0899"
$ws.Range("E11").Value = "  +5.09%  "

# Row 12
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
Set-TextValue $ws "D13" "2.574.72"
$ws.Range("E13").Value = "  +2.12%  "

# Row 14
Set-TextValue $ws "D14" "15.73"
$ws.Range("E14").Value = "  -2.05%  "

# Row 15
Set-TextValue $ws "D15" "22.11"
$ws.Range("E15").Value = "  +0.23%  "

# Row 16
Set-TextValue $ws "D16" "0.807"
$ws.Range("E16").Value = "  -1.59%  "

# Row 17
Set-TextValue $ws "D17" "5.62"
$ws.Range("E17").Value = "  +0.93%  "

# Row 18
Set-TextValue $ws "D18" "2.259.75"
$ws.Range("E18").Value = "  +2.76%  "

# Row 19
Set-TextValue $ws "D19" "42.179.67"
$ws.Range("E19").Value = "  +5.85%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws "D20" "0.0₃0907"
$ws.Range("E20").Value = "  -1.36%  "

# Row 21
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws "D21" "72.31"
$ws.Range("E21").Value = "  +0.36%  "

# Row 22
Set-TextValue $ws "D22" "6.05"
$ws.Range("E22").Value = "  -0.03%  "

# Row 23
Set-TextValue $ws "D23" "252.13"
$ws.Range("E23").Value = "  +9.00%  "

# Row 24
$ws.Range("E24").Value = "  -0.07%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws "D25" "2.38"
$ws.Range("E25").Value = "  -0.66%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws "D26" "2.40"
$ws.Range("E26").Value = "  +1.79%  "

# Row 27
Set-TextValue $ws "D27" "9.73"
$ws.Range("E27").Value = "  +0.87%  "

# Row 28
Set-TextValue $ws "D28" "0.144"
$ws.Range("E28").Value = "  +3.14%  "

# Row 29
Set-TextValue $ws "D29" "169.23"
$ws.Range("E29").Value = "  -1.18%  "

# Row 30
Set-TextValue $ws "D30" "20.11"
$ws.Range("E30").Value = "  +0.36%  "

# Row 31
Set-TextValue $ws "D31" "1.42"
$ws.Range("E31").Value = "  -3.12%  "

# Row 32
Set-TextValue $ws "D32" "2.72"
$ws.Range("E32").Value = "  +0.65%  "

# Row 33
$ws.Range("E33").Value = "  -0.09%  "

# Row 34
Set-TextValue $ws "D34" "5.05"
$ws.Range("E34").Value = "  +7.34%  "

# Row 35
Set-TextValue $ws "D35" "4.69"
$ws.Range("E35").Value = "  +2.92%  "

# Row 36
Set-TextValue $ws "D36" "0.0640"
$ws.Range("E36").Value = "  +2.74%  "

# Row 37
Set-TextValue $ws "D37" "6.70"
$ws.Range("E37").Value = "  -4.71%  "

# Row 38
Set-TextValue $ws "D38" "3.72"
$ws.Range("E38").Value = "  -4.19%  "

# Row 39
$ws.Range("E39").Value = "  -3.15%  "

# Row 40
Set-TextValue $ws "D40" "0.000272"
$ws.Range("E40").Value = "  +41.33%  "

# Row 41
Set-TextValue $ws "D41" "0.998"
$ws.Range("E41").Value = "  -0.13%  "

# Row 42
Set-TextValue $ws "D42" "0.0242"
$ws.Range("E42").Value = "  +5.55%  "

# Row 43
Set-TextValue $ws "D43" "4.85"
$ws.Range("E43").Value = "  -2.93%  "

# Row 44
Set-TextValue $ws "D44" "8.55"
$ws.Range("E44").Value = "  +7.22%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws "D45" "1.23"
$ws.Range("E45").Value = "  +0.04%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D46" "99.49"
$ws.Range("E46").Value = "  -3.70%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D47" "0.0972"
$ws.Range("E47").Value = "  +5.06%  "

# Row 48
Set-TextValue $ws "D48" "1.480.84"
$ws.Range("E48").Value = "  -2.33%  "

# Row 49
Set-TextValue $ws "D49" "16.55"
$ws.Range("E49").Value = "  -7.49%  "

# Row 50
$ws.Range("E50").Value = "  +0.08%  "

# Row 51
$ws.Range("E51").Value = "  +6.21%  "
